# Generate Report for Handoff
#
# Inserts a new "b7327a2c-8ede-44e8-98b2-0a6164f9e85e" file entry ahead of the
# existing "ddc650a5-6cb0-4195-b437-e4d2e34184a2" row on every sheet
# (Overview, zh-cn, de-de), pushing the ".localization-config" row down by
# one, and fills in the new row's handoff data + hyperlinks.

$wb = $excel.ActiveWorkbook

$newFile   = "b7327a2c-8ede-44e8-98b2-0a6164f9e85e.md"
$newXlfBase = "b7327a2c-8ede-44e8-98b2-0a6164f9e85e.3e1266ef0cfcfc25718eed81467f149e69465adc"

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(8).Insert()

$wsOverview.Range("A8").Value = $newFile
$wsOverview.Range("B8").Value = "Ready for handoff"
$wsOverview.Range("C8").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3e1266ef0cfcfc25718eed81467f149e69465adc/e2e/$newFile",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile
)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows.Item(8).Insert()

$zhXlf = "$newXlfBase.zh-cn.xlf"

$wsZh.Range("A8").Value = $newFile
$wsZh.Range("B8").Value = "Ready for handoff"
$wsZh.Range("C8").Value = $zhXlf
$wsZh.Range("D8").Value = "2016-03-09 08:14:58"
$wsZh.Range("G8").Value = "0001-01-01 00:00:00"
$wsZh.Range("H8").Value = "Include"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3e1266ef0cfcfc25718eed81467f149e69465adc/e2e/$newFile",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C8"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e1266ef0cfcfc25718eed81467f149e69465adc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $zhXlf
)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows.Item(8).Insert()

$deXlf = "$newXlfBase.de-de.xlf"

$wsDe.Range("A8").Value = $newFile
$wsDe.Range("B8").Value = "Ready for handoff"
$wsDe.Range("C8").Value = $deXlf
$wsDe.Range("D8").Value = "2016-03-09 08:15:03"
$wsDe.Range("G8").Value = "0001-01-01 00:00:00"
$wsDe.Range("H8").Value = "Include"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3e1266ef0cfcfc25718eed81467f149e69465adc/e2e/$newFile",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C8"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e1266ef0cfcfc25718eed81467f149e69465adc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $deXlf
)
